$d = $word.ActiveDocument

# The first paragraph in the document ("Organisation") is the anchor we
# insert the new heading paragraph before.
$firstPara = $d.Paragraphs(1)

# Insert a brand new paragraph immediately before the first paragraph.
# Word automatically clones the paragraph formatting (pPr/rPr) of the
# paragraph it is inserted before, which already matches the heading
# style used elsewhere in this document (Segoe UI, bold, kern 36, etc.).
$firstPara.Range.InsertParagraphBefore()

# Grab the freshly created (now first) paragraph and set its text.
$newPara = $d.Paragraphs(1)
$newRange = $newPara.Range
$newRange.InsertBefore("Biology")

# Bump the heading up to 28pt (56 half-points) and underline it so it
# stands out as the top-level "science name" heading.
$newRange.Font.Size = 28
$newRange.Font.SizeBi = 28
$newRange.Font.Underline = 1
